# Auto-generated edit script applying scheduled-runner cell updates
# to the Lamia_Profits leve-crafting profit sheets (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 89.75
$ws.Range("I4").Value = 89.75
$ws.Range("K4").Value = 89.75
$ws.Range("M4").Value = 24.25

# Row 33
$ws.Range("H33").Value = 184.8
$ws.Range("I33").Value = 198.11111
$ws.Range("K33").Value = 198.11111
$ws.Range("M33").Value = 30.88889

# Row 47
$ws.Range("H47").Value = 267500
$ws.Range("I47").Value = 500000
$ws.Range("K47").Value = 500000
$ws.Range("M47").Value = -499028

# Row 100
$ws.Range("H100").Value = 6908.4443
$ws.Range("J100").Value = 7744.25
$ws.Range("L100").Value = 7744.25
$ws.Range("N100").Value = -8826.25

# Row 138
$ws.Range("H138").Value = 4029.0588
$ws.Range("J138").Value = 3945.4546
$ws.Range("L138").Value = 11836.3638
$ws.Range("N138").Value = -22116.3638

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 1585.4478
$ws.Range("I32").Value = 1216.8413
$ws.Range("K32").Value = 1216.8413
$ws.Range("M32").Value = -929.8413

# Row 55
$ws.Range("H55").Value = 18038.834
$ws.Range("I55").Value = 3656
$ws.Range("K55").Value = 3656
$ws.Range("M55").Value = -3341

# Row 61
$ws.Range("H61").Value = 5061.647
$ws.Range("I61").Value = 5164.4194
$ws.Range("K61").Value = 5164.4194
$ws.Range("M61").Value = -4952.4194

# Row 74
$ws.Range("H74").Value = 27780174
$ws.Range("J74").Value = 2000
$ws.Range("L74").Value = 2000
$ws.Range("N74").Value = -3748

# Row 77
$ws.Range("H77").Value = 27780174
$ws.Range("J77").Value = 2000
$ws.Range("L77").Value = 10000
$ws.Range("N77").Value = -18736

# Row 132
$ws.Range("H132").Value = 2215.6094
$ws.Range("I132").Value = 1512.2456
$ws.Range("K132").Value = 4536.7368
$ws.Range("M132").Value = -2006.7368

# Row 136
$ws.Range("H136").Value = 5061.647
$ws.Range("I136").Value = 5164.4194
$ws.Range("K136").Value = 15493.2582
$ws.Range("M136").Value = -12943.2582

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 2419.3547
$ws.Range("I20").Value = 2147.2942
$ws.Range("J20").Value = 2749.7144
$ws.Range("K20").Value = 2147.2942
$ws.Range("L20").Value = 2749.7144
$ws.Range("M20").Value = -1900.2942
$ws.Range("N20").Value = -3243.7144

# Row 86
$ws.Range("H86").Value = 1612.5
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").ClearContents()

# Row 89
$ws.Range("H89").Value = 1612.5
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").ClearContents()

# Row 134
$ws.Range("H134").Value = 2788.1538
$ws.Range("I134").Value = 1281.6522
$ws.Range("K134").Value = 3844.9566
$ws.Range("M134").Value = -1309.9566

$ws = $wb.Worksheets.Item("CRP")
# Row 124
$ws.Range("H124").Value = 32500
$ws.Range("I124").Value = 15000
$ws.Range("J124").Value = 38333.332
$ws.Range("K124").Value = 15000
$ws.Range("L124").Value = 38333.332
$ws.Range("M124").Value = -12545
$ws.Range("N124").Value = -43243.332

# Row 132
$ws.Range("H132").Value = 2616.7222
$ws.Range("I132").Value = 1473.1333
$ws.Range("K132").Value = 4419.3999
$ws.Range("M132").Value = -1889.3999

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 1015.9167
$ws.Range("I23").Value = 137.75
$ws.Range("J23").Value = 1455
$ws.Range("K23").Value = 413.25
$ws.Range("L23").Value = 4365
$ws.Range("M23").Value = -178.25
$ws.Range("N23").Value = -4835

# Row 26
$ws.Range("H26").Value = 836.375
$ws.Range("I26").Value = 1007
$ws.Range("J26").Value = 617
$ws.Range("K26").Value = 3021
$ws.Range("L26").Value = 1851
$ws.Range("M26").Value = -2733
$ws.Range("N26").Value = -2427

# Row 68
$ws.Range("H68").Value = 1098
$ws.Range("I68").Value = 1478.4
$ws.Range("J68").Value = 622.5
$ws.Range("K68").Value = 4435.200000000001
$ws.Range("L68").Value = 1867.5
$ws.Range("M68").Value = -3624.200000000001
$ws.Range("N68").Value = -3489.5

# Row 71
$ws.Range("H71").Value = 1098
$ws.Range("I71").Value = 1478.4
$ws.Range("J71").Value = 622.5
$ws.Range("K71").Value = 13305.6
$ws.Range("L71").Value = 5602.5
$ws.Range("M71").Value = -9249.6
$ws.Range("N71").Value = -13714.5

# Row 131
$ws.Range("H131").Value = 11016774
$ws.Range("I131").Value = 8929929
$ws.Range("J131").Value = 12964496
$ws.Range("K131").Value = 26789787
$ws.Range("L131").Value = 38893488
$ws.Range("M131").Value = -26784747
$ws.Range("N131").Value = -38903568

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 10626.04
$ws.Range("I70").Value = 7366.2144
$ws.Range("K70").Value = 7366.2144
$ws.Range("M70").Value = -7096.2144

# Row 73
$ws.Range("H73").Value = 10626.04
$ws.Range("I73").Value = 7366.2144
$ws.Range("K73").Value = 7366.2144
$ws.Range("M73").Value = -6430.2144

# Row 97
$ws.Range("H97").Value = 1110.25
$ws.Range("I97").Value = 990.3
$ws.Range("K97").Value = 990.3
$ws.Range("M97").Value = -494.3

# Row 132
$ws.Range("H132").Value = 4347.1113
$ws.Range("I132").Value = 3687.8
$ws.Range("J132").Value = 7643.6665
$ws.Range("K132").Value = 11063.4
$ws.Range("L132").Value = 22930.9995
$ws.Range("M132").Value = -8533.400000000001
$ws.Range("N132").Value = -27990.9995

$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1954.9375
$ws.Range("I16").Value = 798
$ws.Range("K16").Value = 798
$ws.Range("M16").Value = -628

# Row 40
$ws.Range("H40").Value = 8553
$ws.Range("I40").Value = 7481.9473
$ws.Range("K40").Value = 7481.9473
$ws.Range("M40").Value = -7345.9473

# Row 46
$ws.Range("H46").Value = 2476.4614
$ws.Range("J46").Value = 2763.2727
$ws.Range("L46").Value = 2763.2727
$ws.Range("N46").Value = -3139.2727

# Row 55
$ws.Range("H55").Value = 1283233.5
$ws.Range("I55").Value = 2084124.6
$ws.Range("K55").Value = 2084124.6
$ws.Range("M55").Value = -2083951.6

# Row 61
$ws.Range("H61").Value = 9134.5
$ws.Range("J61").Value = 10702.5
$ws.Range("L61").Value = 10702.5
$ws.Range("N61").Value = -11106.5

# Row 113
$ws.Range("H113").Value = 9134.5
$ws.Range("J113").Value = 10702.5
$ws.Range("L113").Value = 10702.5
$ws.Range("N113").Value = -15042.5

$ws = $wb.Worksheets.Item("WVR")
# Row 5
$ws.Range("H5").Value = 16679999
$ws.Range("I5").Value = 19998
$ws.Range("K5").Value = 19998
$ws.Range("M5").Value = -19886

# Row 122
$ws.Range("H122").Value = 10230.8
$ws.Range("I122").Value = 2149.8
$ws.Range("K122").Value = 6449.400000000001
$ws.Range("M122").Value = -3999.400000000001

# Row 123
$ws.Range("H123").Value = 54333.332
$ws.Range("J123").Value = 54333.332
$ws.Range("L123").Value = 54333.332
$ws.Range("N123").Value = -64133.332

# Row 132
$ws.Range("H132").Value = 6309.0786
$ws.Range("I132").Value = 4329.05
$ws.Range("K132").Value = 12987.15
$ws.Range("M132").Value = -10457.15
